# Update "想去人数" (F column) counts per the site's regenerated data snapshot.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 1592
$ws.Range("F3").Value = 3324
$ws.Range("F5").Value = 755
$ws.Range("F6").Value = 2359
$ws.Range("F8").Value = 426
$ws.Range("F10").Value = 149
$ws.Range("F12").Value = 1114
$ws.Range("F13").Value = 465
$ws.Range("F14").Value = 228
$ws.Range("F15").Value = 93
$ws.Range("F16").Value = 276
$ws.Range("F17").Value = 4892
$ws.Range("F19").Value = 1383
$ws.Range("F20").Value = 3599
$ws.Range("F21").Value = 155
$ws.Range("F22").Value = 211
$ws.Range("F23").Value = 3898
$ws.Range("F24").Value = 5256
$ws.Range("F25").Value = 127
$ws.Range("F26").Value = 987
$ws.Range("F27").Value = 577
$ws.Range("F28").Value = 3378
$ws.Range("F29").Value = 393
$ws.Range("F31").Value = 151
$ws.Range("F32").Value = 97
$ws.Range("F33").Value = 902
$ws.Range("F35").Value = 38
$ws.Range("F36").Value = 55
$ws.Range("F37").Value = 1448
$ws.Range("F38").Value = 146
$ws.Range("F39").Value = 1431
$ws.Range("F40").Value = 39
$ws.Range("F41").Value = 926
$ws.Range("F42").Value = 906
$ws.Range("F43").Value = 528
$ws.Range("F45").Value = 2439
$ws.Range("F47").Value = 187
$ws.Range("F48").Value = 375
$ws.Range("F49").Value = 3763

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F6").Value = 1029
$ws.Range("F23").Value = 43

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 2520

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 2520
$ws.Range("F3").Value = 1592
$ws.Range("F4").Value = 3324
$ws.Range("F6").Value = 755
$ws.Range("F8").Value = 2359
$ws.Range("F10").Value = 426
$ws.Range("F12").Value = 1029
$ws.Range("F13").Value = 149
$ws.Range("F15").Value = 1114
$ws.Range("F16").Value = 465
$ws.Range("F17").Value = 228
$ws.Range("F18").Value = 93
$ws.Range("F19").Value = 276
$ws.Range("F20").Value = 4892
$ws.Range("F22").Value = 1383
$ws.Range("F23").Value = 3898
$ws.Range("F24").Value = 5256
$ws.Range("F25").Value = 127
$ws.Range("F26").Value = 987
$ws.Range("F27").Value = 577
$ws.Range("F28").Value = 3378
$ws.Range("F29").Value = 393
$ws.Range("F31").Value = 151
$ws.Range("F32").Value = 97
$ws.Range("F34").Value = 38
$ws.Range("F35").Value = 55
$ws.Range("F36").Value = 1448
$ws.Range("F37").Value = 1431
$ws.Range("F38").Value = 926
$ws.Range("F39").Value = 528
$ws.Range("F42").Value = 43
$ws.Range("F43").Value = 2439
$ws.Range("F46").Value = 187
$ws.Range("F47").Value = 375
$ws.Range("F49").Value = 3763
